# Add "Average" and "Count" summary rows below the existing
# Minimum Expense / Maximum Expense rows (rows 12-13), mirroring their
# layout: a label in column A and a formula across B:F using a shared
# formula for C:F.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14: Average
$ws.Range("A14").Value = "Average"
$ws.Range("B14").Formula = "=AVERAGE(B4:B8)"
$ws.Range("C14:F14").Formula = "=AVERAGE(C4:C8)"

# Row 15: Count
$ws.Range("A15").Value = "Count"
$ws.Range("B15").Formula = "=COUNT(B4:B8)"
$ws.Range("C15:F15").Formula = "=COUNT(C4:C8)"

# Move the active selection, as recorded when the workbook was last saved.
$ws.Range("A17").Select() | Out-Null
